# Se aplica la tabla del historial de las asesorias del estudiante
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("asesorias")

# Columnas existentes (fila 1 encabezado, fila 2 primer registro):
#   A: Estudiante | B: Usuario | C: Asesor | D: Motivo asesoria | E: Fecha | F: Hora

# Nuevo registro 1 (fila 3)
$ws.Cells.Item(3, 1).Value = "Sebastian Palacio"
$ws.Cells.Item(3, 2).Value = "sebasx200"
$ws.Cells.Item(3, 3).Value = "Juan Carlos Gil"
$ws.Cells.Item(3, 4).Value = "Consulta general"
# "02-11-2023" es ambigua para el parser de fechas de Excel (podria
# interpretarse como una fecha valida); se fuerza formato de texto para
# que el valor se conserve tal cual como cadena.
$ws.Cells.Item(3, 5).NumberFormat = "@"
$ws.Cells.Item(3, 5).Value = "02-11-2023"
$ws.Cells.Item(3, 6).Value = "02:20 - 02:40"

# Nuevo registro 2 (fila 4)
$ws.Cells.Item(4, 1).Value = "Sebastian Palacio"
$ws.Cells.Item(4, 2).Value = "sebasx200"
$ws.Cells.Item(4, 3).Value = "Daniel Henao"
$ws.Cells.Item(4, 4).Value = "Asesoría académica"
$ws.Cells.Item(4, 5).Value = "25-11-2023"
$ws.Cells.Item(4, 6).Value = "00:20 - 00:40"
